$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 21: SA10 / November 18, 2024
$ws.Range("A21").Value = "SA10"

# Match formatting of existing deadline cells (column B, style index 1 -> text number format)
$ws.Range("B21").NumberFormat = $ws.Range("B20").NumberFormat
$ws.Range("B21").Value = "November 18, 2024"

# Update selection to match final saved state
$ws.Range("C24").Select()
